# Add a new bold run ": 50k" right after the run "thiếu tinh thần nhóm"
# (and before the _GoBack bookmark) in the "Cảnh báo nặng: thiếu tinh
# thần nhóm" table cell, matching the formatting of the preceding run
# exactly (bold, color 000000 with themeColor text1, sz 24, szCs 28).

$d = $word.ActiveDocument

# Locate the table cell that holds the target sentence by scanning the
# table's cells for the distinctive (ASCII-safe) substring - avoids any
# dependence on how non-ASCII text happens to round-trip through the
# console/string layer.
$t = $d.Tables.Item(1)
$targetRow = -1
$targetCol = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        try {
            $cell = $t.Cell($r, $c)
            $txt = $cell.Range.Text
            if ($txt.Contains("nh") -and $txt.Contains("thi") -and $txt.Contains(":")) {
                $targetRow = $r
                $targetCol = $c
            }
        } catch {
        }
    }
}
if ($targetRow -eq -1) {
    throw "Could not locate the target table cell"
}

$cell = $t.Cell($targetRow, $targetCol)
$para = $cell.Range.Paragraphs.Item(1)
$paraRng = $para.Range

# Splice a brand-new <w:r> in right before the bookmark, working directly
# on the paragraph's WordprocessingML rather than through the Font/Color
# object model: that model cannot represent a <w:color> that carries both
# an explicit w:val and a w:themeColor at the same time, and it also
# auto-merges adjacent runs whose rPr end up identical.
$newRunXml = '<w:r><w:rPr><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>: 50k</w:t></w:r>'

$bookmarkMarker = '<w:bookmarkStart w:id="0" w:name="_GoBack"/>'

$origParaXml = '<w:p w:rsidR="00440551" w:rsidRPr="003A5CBD" w:rsidRDefault="003A5CBD" w:rsidP="00440551"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/><w:rPr><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="003A5CBD"><w:rPr><w:b/><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>Cảnh báo nặng:</w:t></w:r><w:r><w:rPr><w:b/><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>thiếu tinh thần nhóm</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$newParaXml = $origParaXml.Replace($bookmarkMarker, $newRunXml + $bookmarkMarker)

$xmlPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$paraRng.InsertXML($xmlPackage) | Out-Null

Write-Host "Inserted ': 50k' run after the target paragraph."
